# Fruta / hortaliza, semanal
# A new weekly price record (Primera quality, week of 2021-12-09) was added
# at the top of the data table (row 96), pushing all subsequent records
# down by one row (96:221 -> 97:222).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 96, shifting existing rows 96:221 down to 97:222
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with the new record
$ws.Cells.Item(96, 1).Value = 11
$ws.Cells.Item(96, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(96, 3).Value = "Bíobío"
$ws.Cells.Item(96, 4).Value = 44539
$ws.Cells.Item(96, 5).Value = 8
$ws.Cells.Item(96, 6).Value = "Fruta"
$ws.Cells.Item(96, 7).Value = 100101
$ws.Cells.Item(96, 8).Value = "Berries"
$ws.Cells.Item(96, 9).Value = 100112025
$ws.Cells.Item(96, 10).Value = "Frutilla"
$ws.Cells.Item(96, 11).Value = "Sin especificar"
$ws.Cells.Item(96, 12).Value = "Primera"
$ws.Cells.Item(96, 13).Value = 220
$ws.Cells.Item(96, 14).Value = 6000
$ws.Cells.Item(96, 15).Value = 6500
$ws.Cells.Item(96, 16).Value = 6227
$ws.Cells.Item(96, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(96, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(96, 19).Value = 890
$ws.Cells.Item(96, 20).Value = 7
